$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'53.530.85"
$ws.Range("E2").Value = "  +3.90%  "

# Row 3
$ws.Range("D3").Value = "'3.143.28"
$ws.Range("E3").Value = "  +2.84%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'396.61"
$ws.Range("E5").Value = "  +2.73%  "

# Row 6
$ws.Range("D6").Value = "'109.76"
$ws.Range("E6").Value = "  +6.54%  "

# Row 7
$ws.Range("D7").Value = "'0.546"
$ws.Range("E7").Value = "  +0.57%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  +4.57%  "

# Row 10
$ws.Range("D10").Value = "'38.91"
$ws.Range("E10").Value = "  +5.79%  "

# Row 11
$ws.Range("D11").Value = "'0.140"
$ws.Range("E11").Value = "  +1.38%  "

# Row 12
$ws.Range("D12").Value = "'0.0872"
$ws.Range("E12").Value = "  +1.40%  "

# Row 13
$ws.Range("D13").Value = "'3.653.10"
$ws.Range("E13").Value = "  +3.15%  "

# Row 14
$ws.Range("D14").Value = "'19.09"
$ws.Range("E14").Value = "  +2.84%  "

# Row 15
$ws.Range("D15").Value = "'8.01"
$ws.Range("E15").Value = "  +3.10%  "

# Row 16
$ws.Range("E16").Value = "  +8.88%  "

# Row 17
$ws.Range("D17").Value = "'3.151.66"
$ws.Range("E17").Value = "  +2.96%  "

# Row 18
$ws.Range("D18").Value = "'10.52"
$ws.Range("E18").Value = "  -1.19%  "

# Row 19
$ws.Range("D19").Value = "'53.488.78"
$ws.Range("E19").Value = "  +3.68%  "

# Row 20
$ws.Range("D20").Value = "'3.28"
$ws.Range("E20").Value = "  +3.95%  "

# Row 21
$ws.Range("D21").Value = "'12.82"
$ws.Range("E21").Value = "  +3.00%  "

# Row 22
$ws.Range("D22").Value = "'0.0₃0974"
$ws.Range("E22").Value = "  +0.80%  "

# Row 23
$ws.Range("D23").Value = "'70.92"
$ws.Range("E23").Value = "  +1.08%  "

# Row 24
$ws.Range("D24").Value = "'270.94"
$ws.Range("E24").Value = "  +1.10%  "

# Row 25
$ws.Range("E25").Value = "  +3.33%  "

# Row 26
$ws.Range("D26").Value = "'8.03"
$ws.Range("E26").Value = "  -1.79%  "

# Row 27
$ws.Range("D27").Value = "'27.51"
$ws.Range("E27").Value = "  +2.55%  "

# Row 28
$ws.Range("D28").Value = "'7.34"
$ws.Range("E28").Value = "  +0.85%  "

# Row 29
$ws.Range("E29").Value = "  +0.21%  "

# Row 30
$ws.Range("E30").Value = "  -0.11%  "

# Row 31
$ws.Range("E31").Value = "  +2.99%  "

# Row 32
$ws.Range("D32").Value = "'10.99"
$ws.Range("E32").Value = "  +7.01%  "

# Row 33
$ws.Range("D33").Value = "'37.27"
$ws.Range("E33").Value = "  +7.30%  "

# Row 34
$ws.Range("D34").Value = "'0.0499"
$ws.Range("E34").Value = "  +11.53%  "

# Row 35
$ws.Range("E35").Value = "  +0.71%  "

# Row 36
$ws.Range("D36").Value = "'50.43"
$ws.Range("E36").Value = "  +0.82%  "

# Row 37
$ws.Range("D37").Value = "'3.65"
$ws.Range("E37").Value = "  +9.98%  "

# Row 38
$ws.Range("E38").Value = "  -0.11%  "

# Row 39
$ws.Range("E39").Value = "  +8.87%  "

# Row 40
$ws.Range("D40").Value = "'4.11"
$ws.Range("E40").Value = "  +9.53%  "

# Row 41
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.291"
$ws.Range("E41").Value = "  -0.58%  "

# Row 42
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'17.26"
$ws.Range("E42").Value = "  +1.95%  "

# Row 43
$ws.Range("D43").Value = "'1.89"
$ws.Range("E43").Value = "  +1.43%  "

# Row 44
$ws.Range("D44").Value = "'130.32"
$ws.Range("E44").Value = "  +4.08%  "

# Row 45
$ws.Range("E45").Value = "  +1.51%  "

# Row 46
$ws.Range("D46").Value = "'22.19"
$ws.Range("E46").Value = "  +1.29%  "

# Row 47
$ws.Range("E47").Value = "  -0.91%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  -0.98%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'2.080.30"
$ws.Range("E49").Value = "  +2.36%  "

# Row 50
$ws.Range("D50").Value = "'0.0510"
$ws.Range("E50").Value = "  +18.76%  "

# Row 51
$ws.Range("E51").Value = "  +6.31%  "
